$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the column headers: "_old" -> "_FV2404" and "_new" -> "_FV2410".
#    Columns A-J (1-10) carry the "_old"/"_FV2404" suffixed names, column K (11)
#    is the unchanged "diff" header, and columns L-U (12-21) carry the
#    "_new"/"_FV2410" suffixed names.
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2404"
}
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2410"
}

# 2. Turn the used range A1:U57 into an Excel Table ("Table1") with an
#    autofilter and header row, matching the data already present.
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U57"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3. Freeze the header row (split below row 1, keep top row visible while
#    scrolling).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
